# Updated symbol list on Fri Feb  3 21:51:00 UTC 2023 with GitHub Actions
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the crypto rows that moved since the last snapshot. Values must stay
# plain text (not get auto-converted to numbers/percentages by Excel),
# so each cell is temporarily forced to Text format before the value is
# written, then its style is reset back to Normal so no stray number
# format is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($address, $value) {
    $cell = $ws.Range($address)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$updates = @(
    @{ Cell = "D2";  Value = "329.87" },
    @{ Cell = "E2";  Value = "0.44%" },
    @{ Cell = "D3";  Value = "41.13" },
    @{ Cell = "E3";  Value = "1.92%" },
    @{ Cell = "D4";  Value = "5.701" },
    @{ Cell = "E4";  Value = "-2.26%" },
    @{ Cell = "D5";  Value = "0.08058" },
    @{ Cell = "E5";  Value = "-0.20%" },
    @{ Cell = "D6";  Value = "2.034" },
    @{ Cell = "E6";  Value = "4.13%" },
    @{ Cell = "D7";  Value = "8.707" },
    @{ Cell = "E7";  Value = "-0.72%" },
    @{ Cell = "D8";  Value = "4.518" },
    @{ Cell = "E8";  Value = "-1.73%" },
    @{ Cell = "E9";  Value = "0.12%" },
    @{ Cell = "D10"; Value = "0.9243" },
    @{ Cell = "E10"; Value = "-2.00%" },
    @{ Cell = "E11"; Value = "-1.45%" },
    @{ Cell = "D12"; Value = "0.1945" },
    @{ Cell = "E12"; Value = "-1.62%" },
    @{ Cell = "D13"; Value = "8.267" },
    @{ Cell = "E13"; Value = "-7.45%" },
    @{ Cell = "D14"; Value = "0.09408" },
    @{ Cell = "E14"; Value = "2.15%" },
    @{ Cell = "E15"; Value = "5.39%" },
    @{ Cell = "D16"; Value = "0.1054" },
    @{ Cell = "E16"; Value = "9.34%" },
    @{ Cell = "D17"; Value = "0.001305" },
    @{ Cell = "E17"; Value = "-0.93%" },
    @{ Cell = "D18"; Value = "0.006300" },
    @{ Cell = "E18"; Value = "2.27%" },
    @{ Cell = "E19"; Value = "0.50%" },
    @{ Cell = "E20"; Value = "-2.52%" },
    @{ Cell = "D21"; Value = "0.1418" },
    @{ Cell = "E21"; Value = "0.20%" },
    @{ Cell = "E22"; Value = "9.94%" },
    @{ Cell = "D23"; Value = "0.04420" },
    @{ Cell = "E23"; Value = "0.28%" },
    @{ Cell = "D24"; Value = "0.001259" },
    @{ Cell = "E24"; Value = "-0.24%" },
    @{ Cell = "D25"; Value = "0.004390" },
    @{ Cell = "E25"; Value = "1.01%" },
    @{ Cell = "D26"; Value = "0.0001241" },
    @{ Cell = "E26"; Value = "8.52%" },
    @{ Cell = "D39"; Value = "0.02816" },
    @{ Cell = "E39"; Value = "16.57%" },
    @{ Cell = "D40"; Value = "0.05464" },
    @{ Cell = "E40"; Value = "3.47%" },
    @{ Cell = "D41"; Value = "0.007621" },
    @{ Cell = "E41"; Value = "1.90%" },
    @{ Cell = "D42"; Value = "0.009946" },
    @{ Cell = "E42"; Value = "14.21%" },
    @{ Cell = "E43"; Value = "-0.59%" },
    @{ Cell = "D44"; Value = "0.002132" },
    @{ Cell = "E44"; Value = "1.15%" },
    @{ Cell = "D45"; Value = "0.01184" },
    @{ Cell = "E45"; Value = "7.82%" },
    @{ Cell = "D46"; Value = "0.00006720" },
    @{ Cell = "E46"; Value = "-2.64%" },
    @{ Cell = "E47"; Value = "-0.30%" },
    @{ Cell = "D48"; Value = "0.003000" },
    @{ Cell = "E48"; Value = "-5.24%" },
    @{ Cell = "E49"; Value = "33.65%" },
    @{ Cell = "D50"; Value = "0.00002102" },
    @{ Cell = "E50"; Value = "-0.30%" },
    @{ Cell = "D51"; Value = "0.0002002" },
    @{ Cell = "E51"; Value = "-0.30%" }
)

foreach ($u in $updates) {
    Set-TextValue $u.Cell $u.Value
}

Write-Host "Applied $($updates.Count) cell updates to cryptos sheet"
